$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("NSAA", "AD", "overall", 10),
    @("NSAA", "position", "dhc", 60),
    @("allmatfiles", "jointAngle", "dhc", 60),
    @("NSAA", "position", "dhc", 60),
    @("NSAA", "position", "overall", 60),
    @("NSAA", "position", "acts", 60),
    @("NSAA", "sensorMagneticField", "dhc", 60),
    @("NSAA", "sensorMagneticField", "overall", 60),
    @("NSAA", "sensorMagneticField", "acts", 60),
    @("NSAA", "jointAngle", "dhc", 60),
    @("NSAA", "jointAngle", "overall", 60),
    @("NSAA", "jointAngle", "acts", 60),
    @("NSAA", "jointAngleXZY", "dhc", 60),
    @("NSAA", "jointAngleXZY", "overall", 60),
    @("NSAA", "jointAngleXZY", "acts", 60),
    @("NSAA", "AD", "dhc", 10),
    @("NSAA", "AD", "overall", 10),
    @("NSAA", "AD", "acts", 10),
    @("NSAA", "position", "dhc", 60),
    @("NSAA", "position", "overall", 60),
    @("NSAA", "position", "acts", 60),
    @("NSAA", "sensorMagneticField", "dhc", 60),
    @("NSAA", "sensorMagneticField", "overall", 60),
    @("NSAA", "sensorMagneticField", "acts", 60),
    @("NSAA", "jointAngle", "dhc", 60),
    @("NSAA", "jointAngle", "overall", 60),
    @("NSAA", "jointAngle", "acts", 60),
    @("NSAA", "jointAngleXZY", "dhc", 60),
    @("NSAA", "jointAngleXZY", "overall", 60),
    @("NSAA", "jointAngleXZY", "acts", 60),
    @("NSAA", "AD", "dhc", 10),
    @("NSAA", "AD", "overall", 10),
    @("NSAA", "AD", "acts", 10),
    @("allmatfiles", "jointAngle", "dhc", 60),
    @("allmatfiles", "jointAngle", "overall", 60),
    @("allmatfiles", "jointAngle", "acts", 60)
)

$startRow = 241
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
